$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 (item 13 in the "FICHE DE REVUE DE CODE" table) ---
# LIGNE
$ws.Range("B24").Value = 251
# TYPE
$ws.Range("C24").Value = "Défaut"
# DESCRIPTION / COMMENTAIRES (two-colour rich text, like the rows above it)
$ws.Range("D24").Value = "Questions No.251 est une question de (*)    `"251`t *       Mauvaise  `t Question `"  (questions.bd)"
$ws.Range("D24").Characters(83, 14).Font.Color = 255
# SUIVI (qui)
$ws.Range("F24").Value = "Elie"
# SUIVI (état)
$ws.Range("G24").Value = "en cours"
# Highlight column H like the other filled-in rows (yellow fill, no border)
$ws.Range("H24").Interior.Color = 65535

# --- cursor / selection moved while editing ---
$ws.Range("K15").Select()
